$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 417, shifting existing rows 417-443 down to 418-444.
$ws.Rows.Item(417).Insert()

# Populate the newly inserted row 417 with the new record's values.
$ws.Cells.Item(417, 1).Value = 11
$ws.Cells.Item(417, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(417, 3).Value = "Bíobío"
$ws.Cells.Item(417, 4).Value = 44931
$ws.Cells.Item(417, 5).Value = 8
$ws.Cells.Item(417, 6).Value = 100112006
$ws.Cells.Item(417, 7).Value = "Repollo"
$ws.Cells.Item(417, 8).Value = "Crespo record"
$ws.Cells.Item(417, 9).Value = "Primera"
$ws.Cells.Item(417, 10).Value = 1800
$ws.Cells.Item(417, 11).Value = 800
$ws.Cells.Item(417, 12).Value = 900
$ws.Cells.Item(417, 13).Value = 844
$ws.Cells.Item(417, 14).Value = "$/unidad"
$ws.Cells.Item(417, 15).Value = "Región Metropolitana"
$ws.Cells.Item(417, 16).Value = 844
$ws.Cells.Item(417, 17).Value = 1
$ws.Cells.Item(417, 18).Value = "Hortaliza"
